$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.545.85'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '1.919.09'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4893'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2894'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06699'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '110.48'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.07'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.15%  '
$ws.Range("D12").Value = '1.917.16'
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07586'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.276'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6671'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '291.71'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.60%  '
$ws.Range("D17").Value = '30.534.29'
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007562'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.528'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.54%  '
$ws.Range("D22").Value = '2.165.48'
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.426'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.444'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.63'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.087'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.04%  '
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.448'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.89%  '
$ws.Range("E31").Value = '  -1.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.047'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05017'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7389'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.60%  '
$ws.Range("E35").Value = '  -2.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.726'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("E38").Value = '  -2.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.681'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '112.45'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.01%  '
$ws.Range("E41").Value = '  -2.84%  '
$ws.Range("E42").Value = '  +1.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8641'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '70.82'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.832'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.221'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.25'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.108'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1230'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.2518'
$ws.Range("D51").Style = "Normal"
